# Updated panelApp panels to contain time_taken for metadata purposes.
#
# Two logical changes:
#   1. A gene that was missing from the panel export (RNF113A) is inserted
#      as a new data row right before RNF168 (i.e. at worksheet row 41),
#      pushing every row below it down by one.
#   2. A new "time_taken" column (F) is appended with a per-row timestamp
#      used for metadata/telemetry purposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the missing RNF113A row at worksheet row 41 (shifts the old
#    row 41 "RNF168" and everything after it down to row 42.. 57).
# ---------------------------------------------------------------------
$ws.Rows.Item(41).Insert()

# Copy the formatting (borders/alignment/number-format) of the row above
# into the freshly inserted row so it matches the rest of the table.
$ws.Range("A40:E40").Copy()
$ws.Range("A41:E41").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the new row's values.
$ws.Range("B41").Value2 = "RNF113A"
$ws.Range("C41").Value2 = "ring finger protein 113A"
$ws.Range("E41").Value2 = "Chromosome Breakage Disorders"

# geneConfidence (column D) is stored as TEXT ("3") in this sheet, not a
# number, so copy the value from an existing text cell rather than
# assigning a numeric-looking string (which would coerce to a number).
$ws.Range("D40").Copy()
$ws.Range("D41").PasteSpecial(-4163)       # xlPasteValues

# Column A holds the plain 0-based row index (row number - 2). Native row
# insertion only shifts existing cell content down, it does not
# renumber it, so recompute A for row 41 and every row pushed down after
# it (42..57) explicitly.
for ($r = 41; $r -le 57; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 2
}
# Restore column A's header-row-matching numeric style on the new row
# (row insert gives it a slightly different auto style).
$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)       # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Add the "time_taken" column (F) with a header + one timestamp per
#    data row (rows 2..57).
# ---------------------------------------------------------------------
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)        # xlPasteFormats (bold header style)
$ws.Range("F1").Value2 = "time_taken"

$timeTaken = @(
    "2021-10-05 10:50:27.439422",
    "2021-10-05 10:50:27.439435",
    "2021-10-05 10:50:27.439439",
    "2021-10-05 10:50:27.439442",
    "2021-10-05 10:50:27.439445",
    "2021-10-05 10:50:27.439448",
    "2021-10-05 10:50:27.439452",
    "2021-10-05 10:50:27.439455",
    "2021-10-05 10:50:27.439458",
    "2021-10-05 10:50:27.439462",
    "2021-10-05 10:50:27.439465",
    "2021-10-05 10:50:27.439468",
    "2021-10-05 10:50:27.439471",
    "2021-10-05 10:50:27.439474",
    "2021-10-05 10:50:27.439477",
    "2021-10-05 10:50:27.439480",
    "2021-10-05 10:50:27.439483",
    "2021-10-05 10:50:27.439486",
    "2021-10-05 10:50:27.439489",
    "2021-10-05 10:50:27.439492",
    "2021-10-05 10:50:27.439495",
    "2021-10-05 10:50:27.439498",
    "2021-10-05 10:50:27.439501",
    "2021-10-05 10:50:27.439504",
    "2021-10-05 10:50:27.439507",
    "2021-10-05 10:50:27.439511",
    "2021-10-05 10:50:27.439514",
    "2021-10-05 10:50:27.439517",
    "2021-10-05 10:50:27.439520",
    "2021-10-05 10:50:27.439523",
    "2021-10-05 10:50:27.439526",
    "2021-10-05 10:50:27.439529",
    "2021-10-05 10:50:27.439532",
    "2021-10-05 10:50:27.439536",
    "2021-10-05 10:50:27.439539",
    "2021-10-05 10:50:27.439542",
    "2021-10-05 10:50:27.439545",
    "2021-10-05 10:50:27.439548",
    "2021-10-05 10:50:27.439551",
    "2021-10-05 10:50:27.439554",
    "2021-10-05 10:50:27.439557",
    "2021-10-05 10:50:27.439560",
    "2021-10-05 10:50:27.439563",
    "2021-10-05 10:50:27.439566",
    "2021-10-05 10:50:27.439569",
    "2021-10-05 10:50:27.439572",
    "2021-10-05 10:50:27.439575",
    "2021-10-05 10:50:27.439578",
    "2021-10-05 10:50:27.439581",
    "2021-10-05 10:50:27.439584",
    "2021-10-05 10:50:27.439587",
    "2021-10-05 10:50:27.439590",
    "2021-10-05 10:50:27.439594",
    "2021-10-05 10:50:27.439597",
    "2021-10-05 10:50:27.439600",
    "2021-10-05 10:50:27.439603"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value2 = $timeTaken[$i]
}

Write-Host "edit complete"
